# Auto-generated edit script: updates the cryptos list (prices + % change,
# plus a few re-ranked coin rows) to match the target commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / already-safe-as-text values: direct assignment ---
$ws.Cells.Item(2, 4).Value = '57.265.73'
$ws.Cells.Item(2, 5).Value = '  -0.16%  '
$ws.Cells.Item(3, 4).Value = '2.362.26'
$ws.Cells.Item(3, 5).Value = '  +0.10%  '
$ws.Cells.Item(4, 5).Value = '  +0.19%  '
$ws.Cells.Item(5, 5).Value = '  +4.05%  '
$ws.Cells.Item(6, 5).Value = '  -1.29%  '
$ws.Cells.Item(7, 5).Value = '  -0.25%  '
$ws.Cells.Item(8, 5).Value = '  -0.87%  '
$ws.Cells.Item(9, 4).Value = '2.359.23'
$ws.Cells.Item(9, 5).Value = '  -0.81%  '
$ws.Cells.Item(10, 5).Value = '  -0.69%  '
$ws.Cells.Item(11, 5).Value = '  +0.43%  '
$ws.Cells.Item(12, 5).Value = '  -2.08%  '
$ws.Cells.Item(13, 5).Value = '  +3.57%  '
$ws.Cells.Item(14, 4).Value = '2.748.50'
$ws.Cells.Item(14, 5).Value = '  -1.18%  '
$ws.Cells.Item(15, 5).Value = '  -3.95%  '
$ws.Cells.Item(16, 4).Value = '57.454.67'
$ws.Cells.Item(16, 5).Value = '  +0.22%  '
$ws.Cells.Item(17, 5).Value = '  -0.97%  '
$ws.Cells.Item(18, 4).Value = '2.370.77'
$ws.Cells.Item(18, 5).Value = '  -0.59%  '
$ws.Cells.Item(19, 5).Value = '  +3.31%  '
$ws.Cells.Item(20, 5).Value = '  -0.67%  '
$ws.Cells.Item(21, 5).Value = '  -0.52%  '
$ws.Cells.Item(22, 5).Value = '  +1.05%  '
$ws.Cells.Item(23, 5).Value = '  -0.28%  '
$ws.Cells.Item(24, 5).Value = '  +2.00%  '
$ws.Cells.Item(25, 5).Value = '  +1.28%  '
$ws.Cells.Item(26, 5).Value = '  -1.65%  '
$ws.Cells.Item(27, 5).Value = '  -0.38%  '
$ws.Cells.Item(28, 5).Value = '  +2.80%  '
$ws.Cells.Item(29, 5).Value = '  +3.70%  '
$ws.Cells.Item(30, 5).Value = '  +2.85%  '
$ws.Cells.Item(31, 4).Value = '0.0₃0730'
$ws.Cells.Item(31, 5).Value = '  -1.47%  '
$ws.Cells.Item(32, 5).Value = '  -1.87%  '
$ws.Cells.Item(33, 5).Value = '  -0.10%  '
$ws.Cells.Item(34, 5).Value = '  -0.09%  '
$ws.Cells.Item(35, 2).Value = 'SuiNetwork'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(35, 5).Value = '  +4.76%  '
$ws.Cells.Item(36, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(36, 5).Value = '  -0.16%  '
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 5).Value = '  -2.06%  '
$ws.Cells.Item(38, 5).Value = '  +0.39%  '
$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39, 5).Value = '  -0.40%  '
$ws.Cells.Item(40, 2).Value = 'OKB'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(40, 5).Value = '  +0.93%  '
$ws.Cells.Item(41, 2).Value = 'Aave'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(41, 5).Value = '  +0.08%  '
$ws.Cells.Item(42, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(42, 5).Value = '  -2.99%  '
$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(43, 5).Value = '  -0.88%  '
$ws.Cells.Item(44, 2).Value = 'Bittensor'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(44, 5).Value = '  -0.39%  '
$ws.Cells.Item(45, 2).Value = 'Stellar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(45, 5).Value = '  -0.67%  '
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 5).Value = '  +4.95%  '
$ws.Cells.Item(47, 5).Value = '  -0.95%  '
$ws.Cells.Item(48, 2).Value = 'Mantle'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48, 5).Value = '  -0.17%  '
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 5).Value = '  -0.14%  '
$ws.Cells.Item(50, 2).Value = 'Polygon'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(50, 5).Value = '  +8.34%  '
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 5).Value = '  -0.72%  '

# --- Numeric-looking text values (e.g. "1.00", "0.994") must be written
# through a staging cell + PasteSpecial(values-only) so Excel keeps them as
# literal text (same displayed digits, incl. trailing/leading zeros) instead
# of silently re-parsing them into a Double and dropping formatting digits,
# all without touching the destination cell style.
$stage = $ws.Cells.Item(500, 500)
$stage.Formula = '=""&"1.00"'
$stage.Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"537.47"'
$stage.Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"133.87"'
$stage.Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.994"'
$stage.Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.534"'
$stage.Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.102"'
$stage.Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"5.41"'
$stage.Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.355"'
$stage.Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"23.48"'
$stage.Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.0000134"'
$stage.Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"339.62"'
$stage.Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"10.47"'
$stage.Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"4.23"'
$stage.Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"6.78"'
$stage.Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.997"'
$stage.Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"62.35"'
$stage.Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.168"'
$stage.Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"8.70"'
$stage.Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.991"'
$stage.Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"173.24"'
$stage.Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"1.74"'
$stage.Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"6.14"'
$stage.Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"18.56"'
$stage.Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.998"'
$stage.Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.963"'
$stage.Copy()
$ws.Cells.Item(35, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.993"'
$stage.Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"1.26"'
$stage.Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"4.06"'
$stage.Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"1.59"'
$stage.Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"39.33"'
$stage.Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"149.83"'
$stage.Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.375"'
$stage.Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"3.62"'
$stage.Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"284.76"'
$stage.Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.0930"'
$stage.Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"19.03"'
$stage.Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.0503"'
$stage.Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.564"'
$stage.Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.0218"'
$stage.Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"0.382"'
$stage.Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$stage.Formula = '=""&"17.43"'
$stage.Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)
$stage.Clear()
$excel.CutCopyMode = 0

